$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2300
$ws.Range("F3").Value = 377
$ws.Range("F4").Value = 190
$ws.Range("F5").Value = 302
$ws.Range("F6").Value = 478
$ws.Range("F8").Value = 727
$ws.Range("F9").Value = 524
$ws.Range("F10").Value = 700
$ws.Range("F11").Value = 375
$ws.Range("F12").Value = 69
$ws.Range("F14").Value = 9
$ws.Range("F15").Value = 986
$ws.Range("F16").Value = 15833
$ws.Range("G16").Value = "暂时售罄"
$ws.Range("F17").Value = 319
$ws.Range("F18").Value = 30
$ws.Range("F19").Value = 149
$ws.Range("F20").Value = 259
$ws.Range("F21").Value = 156
$ws.Range("F22").Value = 121
$ws.Range("F23").Value = 10
$ws.Range("F24").Value = 132
$ws.Range("F26").Value = 292
$ws.Range("F27").Value = 118

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 68
$ws.Range("F3").Value = 34
$ws.Range("F4").Value = 158
$ws.Range("G4").Value = 380
$ws.Range("F5").Value = 82
$ws.Range("F7").Value = 213
$ws.Range("F8").Value = 3276
$ws.Range("G8").Value = "已售罄"
$ws.Range("F10").Value = 31
$ws.Range("F11").Value = 1
$ws.Range("F14").Value = 114
$ws.Range("F16").Value = 2660

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 258
$ws.Range("F3").Value = 64
$ws.Range("F4").Value = 478
$ws.Range("F5").Value = 186

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 258
$ws.Range("F3").Value = 64
$ws.Range("F4").Value = 68
$ws.Range("F5").Value = 34
$ws.Range("F6").Value = 2300
$ws.Range("F7").Value = 478
$ws.Range("F8").Value = 377
$ws.Range("F9").Value = 190
$ws.Range("F10").Value = 302
$ws.Range("F11").Value = 478
$ws.Range("F12").Value = 158
$ws.Range("G12").Value = 380
$ws.Range("F14").Value = 82
$ws.Range("F16").Value = 186
$ws.Range("F17").Value = 727
$ws.Range("F18").Value = 524
$ws.Range("F19").Value = 700
$ws.Range("F20").Value = 375
$ws.Range("F21").Value = 69
$ws.Range("F23").Value = 9
$ws.Range("F24").Value = 986
$ws.Range("F25").Value = 15837
$ws.Range("G25").Value = "暂时售罄"
$ws.Range("F26").Value = 213
$ws.Range("F27").Value = 3277
$ws.Range("G27").Value = "已售罄"
$ws.Range("F29").Value = 31
$ws.Range("F30").Value = 1
$ws.Range("F31").Value = 319
$ws.Range("F32").Value = 30
$ws.Range("F33").Value = 149
$ws.Range("F36").Value = 259
$ws.Range("F37").Value = 156
$ws.Range("F38").Value = 121
$ws.Range("F39").Value = 10
$ws.Range("F40").Value = 114
$ws.Range("F42").Value = 132
$ws.Range("F44").Value = 292
$ws.Range("F45").Value = 118
$ws.Range("F46").Value = 2660

